$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("AD3").Value = 11
$ws.Range("AM3").Value = 5.5
$ws.Range("G3").Value = 2.75
$ws.Range("I3").Value = 3.5
$ws.Range("N3").Value = 3.75
$ws.Range("AA4").Value = 2.38
$ws.Range("AB4").Value = 1.53
$ws.Range("AF4").Value = 51
$ws.Range("AL4").Value = 101
$ws.Range("AM4").Value = 5
$ws.Range("AP4").Value = 17
$ws.Range("G4").Value = 4.2
$ws.Range("I4").Value = 2.05
$ws.Range("K4").Value = 1.83
$ws.Range("S4").Value = 2.88
$ws.Range("T4").Value = 1.4
$ws.Range("AO5").Value = 15
$ws.Range("G5").Value = 1.91
$ws.Range("H5").Value = 3.25
$ws.Range("I5").Value = 4.5
$ws.Range("L5").Value = 5
$ws.Range("Q5").Value = 1.95
$ws.Range("R5").Value = 1.9
$ws.Range("O6").Value = 1.67
$ws.Range("P6").Value = 2.1
$ws.Range("AA7").Value = 2.25
$ws.Range("AB7").Value = 1.57
$ws.Range("AC7").Value = 5
$ws.Range("AD7").Value = 6.5
$ws.Range("AF7").Value = 11
$ws.Range("AH7").Value = 41
$ws.Range("AJ7").Value = 7
$ws.Range("AM7").Value = 12
$ws.Range("AN7").Value = 29
$ws.Range("AO7").Value = 21
$ws.Range("AP7").Value = 67
$ws.Range("AQ7").Value = 51
$ws.Range("G7").Value = 1.62
$ws.Range("H7").Value = 3.5
$ws.Range("I7").Value = 6
$ws.Range("J7").Value = 2.3
$ws.Range("K7").Value = 2.05
$ws.Range("L7").Value = 6.5
$ws.Range("AA8").Value = 2
$ws.Range("AB8").Value = 1.73
$ws.Range("AE8").Value = 9
$ws.Range("AF8").Value = 15
$ws.Range("AG8").Value = 17
$ws.Range("AI8").Value = 7.5
$ws.Range("AM8").Value = 10
$ws.Range("AP8").Value = 51
$ws.Range("H8").Value = 3.2
$ws.Range("M8").Value = 1.08
$ws.Range("N8").Value = 8
$ws.Range("O8").Value = 1.4
$ws.Range("P8").Value = 2.75
$ws.Range("Q8").Value = 1.78
$ws.Range("R8").Value = 2.1
$ws.Range("S8").Value = 2.3
$ws.Range("T8").Value = 1.6
$ws.Range("W8").Value = 4.33
$ws.Range("X8").Value = 1.2
$ws.Range("Y8").Value = 1.5
$ws.Range("Z8").Value = 2.5
$ws.Range("AC9").Value = 8
$ws.Range("AD9").Value = 17
$ws.Range("AE9").Value = 13
$ws.Range("AG9").Value = 34
$ws.Range("AH9").Value = 41
$ws.Range("AN9").Value = 9.5
$ws.Range("AO9").Value = 10
$ws.Range("AP9").Value = 21
$ws.Range("G9").Value = 3.6
$ws.Range("H9").Value = 3
$ws.Range("I9").Value = 2.2
$ws.Range("J9").Value = 4.33
$ws.Range("K9").Value = 1.91
$ws.Range("L9").Value = 3
$ws.Range("Q9").Value = 2
$ws.Range("R9").Value = 1.85
$ws.Range("S9").Value = 2.6
$ws.Range("T9").Value = 1.48
$ws.Range("U9").Value = 4.2
$ws.Range("V9").Value = 1.21
$ws.Range("W9").Value = 5.5
$ws.Range("X9").Value = 1.14
$ws.Range("AA10").Value = 2.5
$ws.Range("AB10").Value = 1.47
$ws.Range("AC10").Value = 5.8
$ws.Range("AD10").Value = 13
$ws.Range("AE10").Value = 12.5
$ws.Range("AF10").Value = 40
$ws.Range("AG10").Value = 40
$ws.Range("AI10").Value = 4.1
$ws.Range("AJ10").Value = 5.5
$ws.Range("AK10").Value = 23
$ws.Range("AL10").Value = 200
$ws.Range("AN10").Value = 11.75
$ws.Range("AO10").Value = 12.5
$ws.Range("AQ10").Value = 40
$ws.Range("AR10").Value = 80
$ws.Range("G10").Value = 3.05
$ws.Range("H10").Value = 2.52
$ws.Range("I10").Value = 2.85
$ws.Range("J10").Value = 3.8
$ws.Range("L10").Value = 3.75
$ws.Range("N10").Value = 4.1
$ws.Range("O10").Value = 1.8
$ws.Range("P10").Value = 1.91
$ws.Range("S10").Value = 3.35
$ws.Range("T10").Value = 1.28
$ws.Range("W10").Value = 6.2
$ws.Range("AA11").Value = 2.32
$ws.Range("AD11").Value = 14
$ws.Range("AE11").Value = 12.5
$ws.Range("AF11").Value = 45
$ws.Range("AG11").Value = 40
$ws.Range("AH11").Value = 65
$ws.Range("AN11").Value = 11.5
$ws.Range("AO11").Value = 11.5
$ws.Range("AP11").Value = 35
$ws.Range("AQ11").Value = 35
$ws.Range("AR11").Value = 65
$ws.Range("G11").Value = 3.15
$ws.Range("I11").Value = 2.72
$ws.Range("J11").Value = 3.95
$ws.Range("L11").Value = 3.6
$ws.Range("P11").Value = 2.02
$ws.Range("Y11").Value = 1.7
$ws.Range("AC12").Value = 13
$ws.Range("AI12").Value = 17
$ws.Range("I12").Value = 2.05
$ws.Range("N12").Value = 17
$ws.Range("U12").Value = 1.93
$ws.Range("V12").Value = 1.88
$ws.Range("AS16").Value = 151
$ws.Range("G16").Value = 2.55
$ws.Range("I16").Value = 2.63
$ws.Range("J16").Value = 3.1
$ws.Range("L16").Value = 3.2
$ws.Range("M16").Value = 1.05
$ws.Range("N16").Value = 8.5
$ws.Range("AA17").Value = 1.53
$ws.Range("AB17").Value = 2.38
$ws.Range("AC17").Value = 15
$ws.Range("AD17").Value = 21
$ws.Range("AE17").Value = 13
$ws.Range("AI17").Value = 15
$ws.Range("AJ17").Value = 7.5
$ws.Range("AK17").Value = 12
$ws.Range("AL17").Value = 34
$ws.Range("AM17").Value = 10
$ws.Range("AR17").Value = 21
$ws.Range("AS17").Value = 126
$ws.Range("G17").Value = 3.3
$ws.Range("H17").Value = 3.7
$ws.Range("I17").Value = 1.95
$ws.Range("K17").Value = 2.3
$ws.Range("L17").Value = 2.5
$ws.Range("M17").Value = 1.03
$ws.Range("N17").Value = 10
$ws.Range("O17").Value = 1.17
$ws.Range("P17").Value = 4.5
$ws.Range("S17").Value = 1.57
$ws.Range("T17").Value = 2.35
$ws.Range("U17").Value = 1.95
$ws.Range("V17").Value = 1.85
$ws.Range("W17").Value = 2.38
$ws.Range("X17").Value = 1.53
$ws.Range("Y17").Value = 1.3
$ws.Range("Z17").Value = 3.4
$ws.Range("AA18").Value = 1.57
$ws.Range("AB18").Value = 2.25
$ws.Range("AD18").Value = 17
$ws.Range("AF18").Value = 29
$ws.Range("AG18").Value = 21
$ws.Range("AH18").Value = 26
$ws.Range("AI18").Value = 13
$ws.Range("AJ18").Value = 6.5
$ws.Range("AL18").Value = 41
$ws.Range("AM18").Value = 10
$ws.Range("AN18").Value = 13
$ws.Range("AP18").Value = 21
$ws.Range("AQ18").Value = 17
$ws.Range("G18").Value = 2.88
$ws.Range("H18").Value = 3.4
$ws.Range("I18").Value = 2.35
$ws.Range("J18").Value = 3.25
$ws.Range("K18").Value = 2.2
$ws.Range("L18").Value = 2.88
$ws.Range("O18").Value = 1.2
$ws.Range("P18").Value = 4.33
$ws.Range("AC20").Value = 8.5
$ws.Range("AD20").Value = 13
$ws.Range("AM20").Value = 8
$ws.Range("AN20").Value = 12
$ws.Range("G20").Value = 2.88
$ws.Range("I20").Value = 2.45
$ws.Range("L20").Value = 3.2
$ws.Range("O20").Value = 1.36
$ws.Range("P20").Value = 3
$ws.Range("S20").Value = 2.15
$ws.Range("T20").Value = 1.67
$ws.Range("AD21").Value = 11
$ws.Range("AF21").Value = 21
$ws.Range("AI21").Value = 10
$ws.Range("AM21").Value = 10
$ws.Range("AN21").Value = 15
$ws.Range("AO21").Value = 11
$ws.Range("G21").Value = 2.15
$ws.Range("H21").Value = 3.4
$ws.Range("I21").Value = 3.1
$ws.Range("J21").Value = 2.88
$ws.Range("K21").Value = 2.1
$ws.Range("L21").Value = 3.75
$ws.Range("S21").Value = 1.95
$ws.Range("T21").Value = 1.85
$ws.Range("W21").Value = 3.4
$ws.Range("X21").Value = 1.3
$ws.Range("AC22").Value = 11
$ws.Range("AE22").Value = 11
$ws.Range("AF22").Value = 29
$ws.Range("AP22").Value = 23
$ws.Range("G22").Value = 2.8
$ws.Range("H22").Value = 3.3
$ws.Range("I22").Value = 2.45
$ws.Range("J22").Value = 3.4
$ws.Range("L22").Value = 3
$ws.Range("M22").Value = 1.04
$ws.Range("N22").Value = 12
$ws.Range("O22").Value = 1.22
$ws.Range("P22").Value = 4
$ws.Range("W22").Value = 2.75
$ws.Range("X22").Value = 1.4
$ws.Range("AA24").Value = 1.47
$ws.Range("AB24").Value = 2.5
$ws.Range("AC24").Value = 12
$ws.Range("AD24").Value = 10
$ws.Range("AF24").Value = 12
$ws.Range("AG24").Value = 10.25
$ws.Range("AH24").Value = 16.5
$ws.Range("AI24").Value = 10.25
$ws.Range("AJ24").Value = 10.25
$ws.Range("AL24").Value = 37
$ws.Range("AM24").Value = 28
$ws.Range("AN24").Value = 45
$ws.Range("AO24").Value = 18
$ws.Range("AP24").Value = 110
$ws.Range("AQ24").Value = 45
$ws.Range("AS24").Value = 175
$ws.Range("G24").Value = 1.45
$ws.Range("H24").Value = 4.65
$ws.Range("I24").Value = 5.6
$ws.Range("J24").Value = 1.9
$ws.Range("K24").Value = 2.62
$ws.Range("L24").Value = 5
$ws.Range("O24").Value = 1.11
$ws.Range("P24").Value = 5.6
$ws.Range("S24").Value = 1.35
$ws.Range("T24").Value = 2.92
$ws.Range("W24").Value = 1.85
$ws.Range("X24").Value = 1.85
$ws.Range("Y24").Value = 1.22
$ws.Range("Z24").Value = 3.85
